$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.84976866666667
$ws.Range("H2").Value = 65.549306
$ws.Range("I2").Value = 0.05020018890879543
$ws.Range("J2").Value = 0.05020018890879543
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 261.380203
$ws.Range("N2").Value = 784.1406089999999
$ws.Range("O2").Value = 0.6968677182772199
$ws.Range("P2").Value = 0.6968677182772199
$ws.Range("Q2").Value = 5711.096969596372
$ws.Range("R2").Value = 51399.87272636735
$ws.Range("S2").Value = 0.03498289110195767
$ws.Range("T2").Value = 0.03498289110195767

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.84976866666667
$ws.Range("H3").Value = 65.549306
$ws.Range("I3").Value = 0.05020018890879543
$ws.Range("J3").Value = 0.05020018890879543
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 31.999428
$ws.Range("N3").Value = 95.998284
$ws.Range("O3").Value = 0.08531391482826334
$ws.Range("P3").Value = 0.08531391482826335
$ws.Range("Q3").Value = 699.1800992656559
$ws.Range("R3").Value = 6292.620893390904
$ws.Range("S3").Value = 0.004282774640927703
$ws.Range("T3").Value = 0.004282774640927704

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.84976866666667
$ws.Range("H4").Value = 65.549306
$ws.Range("I4").Value = 0.05020018890879543
$ws.Range("J4").Value = 0.05020018890879543
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 81.699019
$ws.Range("N4").Value = 245.097057
$ws.Range("O4").Value = 0.2178183668945166
$ws.Range("P4").Value = 0.2178183668945167
$ws.Range("Q4").Value = 1785.104665443605
$ws.Range("R4").Value = 16065.94198899244
$ws.Range("S4").Value = 0.01093452316591005
$ws.Range("T4").Value = 0.01093452316591005

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 385.0524703333334
$ws.Range("H5").Value = 1155.157411
$ws.Range("I5").Value = 0.8846641374295412
$ws.Range("J5").Value = 0.8846641374295412
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 261.380203
$ws.Range("N5").Value = 784.1406089999999
$ws.Range("O5").Value = 0.6968677182772199
$ws.Range("P5").Value = 0.6968677182772199
$ws.Range("Q5").Value = 100645.0928613781
$ws.Range("R5").Value = 905805.8357524034
$ws.Range("S5").Value = 0.6164938788922093
$ws.Range("T5").Value = 0.6164938788922093

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 385.0524703333334
$ws.Range("H6").Value = 1155.157411
$ws.Range("I6").Value = 0.8846641374295412
$ws.Range("J6").Value = 0.8846641374295412
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.999428
$ws.Range("N6").Value = 95.998284
$ws.Range("O6").Value = 0.08531391482826334
$ws.Range("P6").Value = 0.08531391482826335
$ws.Range("Q6").Value = 12321.45880065364
$ws.Range("R6").Value = 110893.1292058827
$ws.Range("S6").Value = 0.07547416087228292
$ws.Range("T6").Value = 0.07547416087228294

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 385.0524703333334
$ws.Range("H7").Value = 1155.157411
$ws.Range("I7").Value = 0.8846641374295412
$ws.Range("J7").Value = 0.8846641374295412
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 81.699019
$ws.Range("N7").Value = 245.097057
$ws.Range("O7").Value = 0.2178183668945166
$ws.Range("P7").Value = 0.2178183668945167
$ws.Range("Q7").Value = 31458.40908975994
$ws.Range("R7").Value = 283125.6818078395
$ws.Range("S7").Value = 0.1926960976650489
$ws.Range("T7").Value = 0.1926960976650489

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.350479
$ws.Range("H8").Value = 85.05143699999999
$ws.Range("I8").Value = 0.06513567366166337
$ws.Range("J8").Value = 0.06513567366166337
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 261.380203
$ws.Range("N8").Value = 784.1406089999999
$ws.Range("O8").Value = 0.6968677182772199
$ws.Range("P8").Value = 0.6968677182772199
$ws.Range("Q8").Value = 7410.253956167236
$ws.Range("R8").Value = 66692.28560550512
$ws.Range("S8").Value = 0.04539094828305296
$ws.Range("T8").Value = 0.04539094828305296

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.350479
$ws.Range("H9").Value = 85.05143699999999
$ws.Range("I9").Value = 0.06513567366166337
$ws.Range("J9").Value = 0.06513567366166337
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 31.999428
$ws.Range("N9").Value = 95.998284
$ws.Range("O9").Value = 0.08531391482826334
$ws.Range("P9").Value = 0.08531391482826335
$ws.Range("Q9").Value = 907.1991115260118
$ws.Range("R9").Value = 8164.792003734107
$ws.Range("S9").Value = 0.005556979315052704
$ws.Range("T9").Value = 0.005556979315052705

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.350479
$ws.Range("H10").Value = 85.05143699999999
$ws.Range("I10").Value = 0.06513567366166337
$ws.Range("J10").Value = 0.06513567366166337
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 81.699019
$ws.Range("N10").Value = 245.097057
$ws.Range("O10").Value = 0.2178183668945166
$ws.Range("P10").Value = 0.2178183668945167
$ws.Range("Q10").Value = 2316.206322480101
$ws.Range("R10").Value = 20845.85690232091
$ws.Range("S10").Value = 0.01418774606355769
$ws.Range("T10").Value = 0.0141877460635577
